# Update the "Mass (KG)" column (E2:E9) from whole-number grams to
# kilograms expressed as a fraction (divide by 1000).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.595
$ws.Range("E3").Value = 0.678
$ws.Range("E4").Value = 0.704
$ws.Range("E5").Value = 0.844
$ws.Range("E6").Value = 0.461
$ws.Range("E7").Value = 0.589
$ws.Range("E8").Value = 0.533
$ws.Range("E9").Value = 0.406
